# Implementing pipelines for interacting with alert classification files
#
# This script edits the "Runs" sheet (2nd sheet) of the workbook:
#  - Inserts two new columns (Precision, Sensitivity) right after the FN column
#  - Fills the new columns with 0 for every model row
#  - Converts the "Active" boolean literal into a =TRUE() formula
#  - Re-creates the hyperlinks on the (now shifted) Documentation column, restoring
#    their original (blue, non-underlined) look
#  - Normalizes the smart-quotes in the custom parameters JSON text
#  - Makes the "Runs" sheet the active / selected sheet of the workbook

$wb = $excel.ActiveWorkbook
$wsAlerts = $wb.Worksheets.Item(1)
$wsRuns = $wb.Worksheets.Item(2)

# Insert 2 new columns (I and J) on the Runs sheet, before the old "F1Score" column
$insertRange = $wsRuns.Range($wsRuns.Cells.Item(1, 9), $wsRuns.Cells.Item(1, 10))
$insertRange.EntireColumn.Insert()

# After the insert, the Documentation (hyperlink) column formatting already
# landed on its new home (column N) - remember that look before (re)adding
# the hyperlinks, which would otherwise reset it to the default style
$wsRuns.Range("N4").Copy() | Out-Null

# New headers for the inserted columns
$wsRuns.Range("I3").Value = "Precision"
$wsRuns.Range("J3").Value = "Sensitivity"

# New data values (0) for the 3 model rows
$wsRuns.Range("I4").Value = 0
$wsRuns.Range("J4").Value = 0
$wsRuns.Range("I5").Value = 0
$wsRuns.Range("J5").Value = 0
$wsRuns.Range("I6").Value = 0
$wsRuns.Range("J6").Value = 0

# The old "Active" boolean literal (now in column M) becomes a =TRUE() formula
$wsRuns.Range("M4").Formula = "=TRUE()"
$wsRuns.Range("M5").Formula = "=TRUE()"
$wsRuns.Range("M6").Formula = "=TRUE()"

# Fix the custom parameters JSON text that used smart quotes
$wsRuns.Range("O4").Value = '{"regParam":1, "tol":0.1}'
$wsRuns.Range("O6").Value = '{"regParam":1, "tol":0.1}'

# The hyperlinks did not automatically follow the column insertion, recreate them
# on their new (shifted) home column N
$wsRuns.Range("L4").Hyperlinks.Delete()
$wsRuns.Hyperlinks.Add($wsRuns.Range("N4"), "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/LinearSVC.html", "", "", "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/LinearSVC.html")
$wsRuns.Hyperlinks.Add($wsRuns.Range("N5"), "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/GBTClassifier.html", "", "", "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/GBTClassifier.html")
$wsRuns.Hyperlinks.Add($wsRuns.Range("N6"), "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/LogisticRegression.html", "", "", "https://spark.apache.org/docs/3.0.1/api/scala/org/apache/spark/ml/classification/LogisticRegression.html")

# Restore the original (blue, non-underlined) hyperlink cell look that Excel
# replaced with its default "Hyperlink" style when the links were (re)created
$wsRuns.Range("N4").PasteSpecial(-4122) | Out-Null
$wsRuns.Range("N5").PasteSpecial(-4122) | Out-Null
$wsRuns.Range("N6").PasteSpecial(-4122) | Out-Null

# Make the Runs sheet the active sheet with its own selection
$wsRuns.Activate() | Out-Null
$wsRuns.Range("G13").Select() | Out-Null
